# Add a new "transit_lines" worksheet, positioned right after "vehicle_types"
# (i.e. before "shunting_locations_on_route"), with a header row describing a
# transit-line based filter/vehicle-allocation strategy.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the "vehicle_types" sheet.
$vehicleTypesSheet = $wb.Worksheets.Item("vehicle_types")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $vehicleTypesSheet)
$newSheet.Name = "transit_lines"

# Header row content.
$newSheet.Range("A1").Value = "transit_line_id"
$newSheet.Range("B1").Value = "vehicle_type"
$newSheet.Range("G1").Value = "Note: This sheet is optional; if list is empty, no filter is applied."

# Headers / note use the same bold style as the other sheets in this workbook.
$newSheet.Range("A1:B1").Font.Bold = $true
$newSheet.Range("G1").Font.Bold = $true

# Match the column sizing used for the other "id" columns in this workbook.
$newSheet.Columns.Item(1).ColumnWidth = 12.45

# Leave the sheet on the selection/active state it was saved with.
$newSheet.Range("F13").Select() | Out-Null
$newSheet.Activate() | Out-Null
